# Money Tracker deploy v1 - with swagger
# Insert a new "cisNumber" column right after UserDetailId (new column B),
# pushing UserName + the boolean flag columns one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column B (UserName).
$ws.Columns.Item(2).Insert()

# Mark the new column as text-formatted (matches the "@" / numFmtId 49 style
# used so the CIS numbers are stored/displayed as text, not numbers).
$ws.Range("B1:B8").NumberFormat = "@"
$ws.Columns.Item(2).ColumnWidth = 13.109375

# Header + CIS number values for the 7 seeded users.
$ws.Range("B1").Value = "cisNumber"
$ws.Range("B2").Value = "110001263706"
$ws.Range("B3").Value = "110001263707"
$ws.Range("B4").Value = "110001263708"
$ws.Range("B5").Value = "110001263709"
$ws.Range("B6").Value = "110001263710"
$ws.Range("B7").Value = "110001263711"
$ws.Range("B8").Value = "110001263712"

# Match the saved selection from the authored workbook.
$ws.Range("B7").Select() | Out-Null
